# update p2map for correct auto names
#
# Corrects the "RealAutoParent" (column B) values on the panel2_v2 sheet so
# each row's auto-generated parent name lines up with the right node, plus
# the follow-on view/selection + column-width changes that came with it.

$wb = $excel.ActiveWorkbook

# --- panel2_v2: fix the RealAutoParent (column B) values, rows 3-15 -------
$ws2v2 = $wb.Worksheets.Item("panel2_v2")

$ws2v2.Range("B3").Value  = "PE-A"
$ws2v2.Range("B4").Value  = "Singlets"
$ws2v2.Range("B5").Value  = "PBMC"
$ws2v2.Range("B6").Value  = "D_NK_M"
$ws2v2.Range("B7").Value  = "CD14+"
$ws2v2.Range("B8").Value  = "CD14+"
$ws2v2.Range("B9").Value  = "D_NK_M"
$ws2v2.Range("B10").Value = "CD20-"
$ws2v2.Range("B11").Value = "Dendritic"
$ws2v2.Range("B12").Value = "Dendritic"
$ws2v2.Range("B13").Value = "D_NK_M"
$ws2v2.Range("B14").Value = "CD20-/CD16+"
$ws2v2.Range("B15").Value = "CD20-/CD16+"

# Column B is now full of much shorter strings than before, so it gets
# resized (auto-fit to its new widest entry, the "RealAutoParent" header).
$ws2v2.Columns.Item(2).ColumnWidth = 32

# --- panel2: move the saved selection to A10 ------------------------------
$ws2 = $wb.Worksheets.Item("panel2")
$ws2.Activate()
$ws2.Range("A10").Select()

# --- panel2_v2: move the saved selection to B14, re-activate as the -------
# tab that was selected when the workbook was saved
$ws2v2.Activate()
$ws2v2.Range("B14").Select()
